{"js": "// Replace each two-digit multiplication expression in the table with its\n// new value, in document order, while preserving run formatting\n// (rFonts TimeNewRoman, sz 30) by replacing the matched range's text\n// in place instead of clearing/re-inserting the whole cell.\nconst oldVals = [\"74\u00d777=\", \"92\u00d731=\", \"13\u00d728=\", \"77\u00d742=\", \"71\u00d716=\", \"14\u00d746=\", \"18\u00d720=\", \"67\u00d745=\", \"79\u00d732=\", \"52\u00d744=\", \"20\u00d796=\", \"88\u00d734=\", \"28\u00d745=\", \"90\u00d775=\", \"21\u00d737=\", \"35\u00d765=\", \"96\u00d762=\", \"91\u00d792=\", \"44\u00d726=\", \"97\u00d749=\", \"29\u00d775=\", \"29\u00d758=\", \"68\u00d791=\", \"15\u00d728=\", \"96\u00d788=\", \"54\u00d716=\", \"74\u00d766=\", \"73\u00d794=\", \"46\u00d744=\", \"32\u00d790=\", \"96\u00d733=\", \"40\u00d745=\", \"37\u00d741=\", \"36\u00d762=\", \"32\u00d776=\", \"96\u00d796=\", \"21\u00d729=\", \"32\u00d751=\", \"78\u00d754=\", \"78\u00d717=\", \"87\u00d713=\", \"68\u00d719=\", \"92\u00d752=\", \"25\u00d770=\", \"58\u00d776=\", \"33\u00d717=\", \"37\u00d791=\", \"43\u00d730=\", \"79\u00d756=\", \"51\u00d794=\", \"39\u00d760=\", \"10\u00d728=\", \"33\u00d799=\", \"39\u00d747=\", \"43\u00d741=\", \"17\u00d724=\", \"87\u00d752=\", \"88\u00d779=\", \"85\u00d757=\", \"32\u00d780=\", \"48\u00d711=\", \"18\u00d761=\", \"43\u00d766=\", \"74\u00d742=\", \"60\u00d769=\", \"100\u00d747=\", \"57\u00d754=\", \"26\u00d792=\", \"76\u00d798=\", \"95\u00d736=\", \"55\u00d723=\", \"82\u00d715=\", \"16\u00d776=\", \"53\u00d717=\", \"91\u00d749=\", \"87\u00d712=\", \"22\u00d799=\", \"97\u00d723=\", \"12\u00d755=\", \"75\u00d757=\", \"47\u00d735=\", \"68\u00d746=\", \"53\u00d740=\", \"19\u00d796=\", \"29\u00d731=\", \"61\u00d766=\", \"28\u00d724=\", \"74\u00d774=\", \"82\u00d794=\", \"74\u00d711=\", \"70\u00d764=\", \"24\u00d723=\", \"10\u00d748=\", \"76\u00d763=\", \"76\u00d772=\", \"38\u00d783=\", \"59\u00d763=\", \"54\u00d738=\", \"99\u00d719=\", \"10\u00d718=\"];\nconst newVals = [\"21\u00d767=\", \"81\u00d781=\", \"41\u00d728=\", \"92\u00d799=\", \"94\u00d780=\", \"62\u00d777=\", \"26\u00d757=\", \"72\u00d733=\", \"56\u00d710=\", \"48\u00d759=\", \"77\u00d748=\", \"79\u00d776=\", \"55\u00d782=\", \"84\u00d760=\", \"32\u00d757=\", \"34\u00d723=\", \"64\u00d734=\", \"90\u00d739=\", \"86\u00d741=\", \"56\u00d781=\", \"91\u00d716=\", \"31\u00d767=\", \"58\u00d750=\", \"73\u00d753=\", \"98\u00d790=\", \"59\u00d788=\", \"96\u00d770=\", \"21\u00d756=\", \"59\u00d751=\", \"76\u00d797=\", \"31\u00d787=\", \"14\u00d739=\", \"17\u00d756=\", \"24\u00d776=\", \"32\u00d712=\", \"60\u00d751=\", \"37\u00d722=\", \"93\u00d751=\", \"64\u00d729=\", \"70\u00d786=\", \"92\u00d758=\", \"68\u00d753=\", \"93\u00d733=\", \"47\u00d737=\", \"40\u00d796=\", \"69\u00d788=\", \"97\u00d776=\", \"57\u00d748=\", \"10\u00d723=\", \"74\u00d778=\", \"53\u00d790=\", \"33\u00d758=\", \"78\u00d722=\", \"72\u00d762=\", \"97\u00d717=\", \"83\u00d746=\", \"78\u00d782=\", \"58\u00d773=\", \"69\u00d759=\", \"21\u00d731=\", \"64\u00d743=\", \"24\u00d774=\", \"65\u00d767=\", \"21\u00d759=\", \"44\u00d793=\", \"18\u00d7100=\", \"58\u00d765=\", \"16\u00d790=\", \"51\u00d723=\", \"89\u00d715=\", \"95\u00d782=\", \"31\u00d717=\", \"24\u00d711=\", \"14\u00d780=\", \"34\u00d727=\", \"44\u00d791=\", \"47\u00d726=\", \"82\u00d788=\", \"96\u00d738=\", \"36\u00d754=\", \"35\u00d723=\", \"16\u00d759=\", \"46\u00d797=\", \"39\u00d717=\", \"19\u00d782=\", \"78\u00d781=\", \"76\u00d712=\", \"31\u00d775=\", \"63\u00d753=\", \"46\u00d753=\", \"99\u00d786=\", \"58\u00d764=\", \"97\u00d798=\", \"43\u00d772=\", \"100\u00d751=\", \"43\u00d762=\", \"16\u00d729=\", \"96\u00d719=\", \"84\u00d727=\", \"72\u00d717=\"];\n\nconst body = context.document.body;\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"No table found in document\");\n}\n\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\n// Collect all cells in row-major (reading) order.\nconst cellList = [];\nfor (const row of rows.items) {\n  const cells = row.cells;\n  cells.load(\"items\");\n  await context.sync();\n  for (const cell of cells.items) {\n    cellList.push(cell);\n  }\n}\n\nif (cellList.length !== oldVals.length) {\n  throw new Error(`Cell count ${cellList.length} does not match expected ${oldVals.length}`);\n}\n\n// Load each cell's body text so we can confirm we are editing the\n// expected cell before mutating it.\nfor (const cell of cellList) {\n  cell.body.load(\"text\");\n}\nawait context.sync();\n\nfor (let i = 0; i < cellList.length; i++) {\n  const cell = cellList[i];\n  const expectedOld = oldVals[i];\n  const newVal = newVals[i];\n  const actual = cell.body.text.trim();\n\n  if (actual !== expectedOld) {\n    throw new Error(`Cell ${i}: expected \"${expectedOld}\" but found \"${actual}\"`);\n  }\n\n  // Search within this cell's body for the exact expression and replace\n  // just that range, which keeps the run's original formatting\n  // (rFonts/sz) intact.\n  const searchResults = cell.body.search(expectedOld, { matchCase: true, matchWholeWord: false });\n  searchResults.load(\"items\");\n  await context.sync();\n\n  if (searchResults.items.length === 0) {\n    throw new Error(`Cell ${i}: could not locate \"${expectedOld}\" via search`);\n  }\n\n  searchResults.items[0].insertText(newVal, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Replace each two-digit multiplication expression in the table with its\n# new value, in document (row-major) order. Assigning to Range.Text only\n# rewrites the run's text content, so the existing run formatting\n# (rFonts TimeNewRoman, sz 30) is preserved.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$oldVals = @(\n    \"74\u00d777=\",\n    \"92\u00d731=\",\n    \"13\u00d728=\",\n    \"77\u00d742=\",\n    \"71\u00d716=\",\n    \"14\u00d746=\",\n    \"18\u00d720=\",\n    \"67\u00d745=\",\n    \"79\u00d732=\",\n    \"52\u00d744=\",\n    \"20\u00d796=\",\n    \"88\u00d734=\",\n    \"28\u00d745=\",\n    \"90\u00d775=\",\n    \"21\u00d737=\",\n    \"35\u00d765=\",\n    \"96\u00d762=\",\n    \"91\u00d792=\",\n    \"44\u00d726=\",\n    \"97\u00d749=\",\n    \"29\u00d775=\",\n    \"29\u00d758=\",\n    \"68\u00d791=\",\n    \"15\u00d728=\",\n    \"96\u00d788=\",\n    \"54\u00d716=\",\n    \"74\u00d766=\",\n    \"73\u00d794=\",\n    \"46\u00d744=\",\n    \"32\u00d790=\",\n    \"96\u00d733=\",\n    \"40\u00d745=\",\n    \"37\u00d741=\",\n    \"36\u00d762=\",\n    \"32\u00d776=\",\n    \"96\u00d796=\",\n    \"21\u00d729=\",\n    \"32\u00d751=\",\n    \"78\u00d754=\",\n    \"78\u00d717=\",\n    \"87\u00d713=\",\n    \"68\u00d719=\",\n    \"92\u00d752=\",\n    \"25\u00d770=\",\n    \"58\u00d776=\",\n    \"33\u00d717=\",\n    \"37\u00d791=\",\n    \"43\u00d730=\",\n    \"79\u00d756=\",\n    \"51\u00d794=\",\n    \"39\u00d760=\",\n    \"10\u00d728=\",\n    \"33\u00d799=\",\n    \"39\u00d747=\",\n    \"43\u00d741=\",\n    \"17\u00d724=\",\n    \"87\u00d752=\",\n    \"88\u00d779=\",\n    \"85\u00d757=\",\n    \"32\u00d780=\",\n    \"48\u00d711=\",\n    \"18\u00d761=\",\n    \"43\u00d766=\",\n    \"74\u00d742=\",\n    \"60\u00d769=\",\n    \"100\u00d747=\",\n    \"57\u00d754=\",\n    \"26\u00d792=\",\n    \"76\u00d798=\",\n    \"95\u00d736=\",\n    \"55\u00d723=\",\n    \"82\u00d715=\",\n    \"16\u00d776=\",\n    \"53\u00d717=\",\n    \"91\u00d749=\",\n    \"87\u00d712=\",\n    \"22\u00d799=\",\n    \"97\u00d723=\",\n    \"12\u00d755=\",\n    \"75\u00d757=\",\n    \"47\u00d735=\",\n    \"68\u00d746=\",\n    \"53\u00d740=\",\n    \"19\u00d796=\",\n    \"29\u00d731=\",\n    \"61\u00d766=\",\n    \"28\u00d724=\",\n    \"74\u00d774=\",\n    \"82\u00d794=\",\n    \"74\u00d711=\",\n    \"70\u00d764=\",\n    \"24\u00d723=\",\n    \"10\u00d748=\",\n    \"76\u00d763=\",\n    \"76\u00d772=\",\n    \"38\u00d783=\",\n    \"59\u00d763=\",\n    \"54\u00d738=\",\n    \"99\u00d719=\",\n    \"10\u00d718=\"\n)\n\n$newVals = @(\n    \"21\u00d767=\",\n    \"81\u00d781=\",\n    \"41\u00d728=\",\n    \"92\u00d799=\",\n    \"94\u00d780=\",\n    \"62\u00d777=\",\n    \"26\u00d757=\",\n    \"72\u00d733=\",\n    \"56\u00d710=\",\n    \"48\u00d759=\",\n    \"77\u00d748=\",\n    \"79\u00d776=\",\n    \"55\u00d782=\",\n    \"84\u00d760=\",\n    \"32\u00d757=\",\n    \"34\u00d723=\",\n    \"64\u00d734=\",\n    \"90\u00d739=\",\n    \"86\u00d741=\",\n    \"56\u00d781=\",\n    \"91\u00d716=\",\n    \"31\u00d767=\",\n    \"58\u00d750=\",\n    \"73\u00d753=\",\n    \"98\u00d790=\",\n    \"59\u00d788=\",\n    \"96\u00d770=\",\n    \"21\u00d756=\",\n    \"59\u00d751=\",\n    \"76\u00d797=\",\n    \"31\u00d787=\",\n    \"14\u00d739=\",\n    \"17\u00d756=\",\n    \"24\u00d776=\",\n    \"32\u00d712=\",\n    \"60\u00d751=\",\n    \"37\u00d722=\",\n    \"93\u00d751=\",\n    \"64\u00d729=\",\n    \"70\u00d786=\",\n    \"92\u00d758=\",\n    \"68\u00d753=\",\n    \"93\u00d733=\",\n    \"47\u00d737=\",\n    \"40\u00d796=\",\n    \"69\u00d788=\",\n    \"97\u00d776=\",\n    \"57\u00d748=\",\n    \"10\u00d723=\",\n    \"74\u00d778=\",\n    \"53\u00d790=\",\n    \"33\u00d758=\",\n    \"78\u00d722=\",\n    \"72\u00d762=\",\n    \"97\u00d717=\",\n    \"83\u00d746=\",\n    \"78\u00d782=\",\n    \"58\u00d773=\",\n    \"69\u00d759=\",\n    \"21\u00d731=\",\n    \"64\u00d743=\",\n    \"24\u00d774=\",\n    \"65\u00d767=\",\n    \"21\u00d759=\",\n    \"44\u00d793=\",\n    \"18\u00d7100=\",\n    \"58\u00d765=\",\n    \"16\u00d790=\",\n    \"51\u00d723=\",\n    \"89\u00d715=\",\n    \"95\u00d782=\",\n    \"31\u00d717=\",\n    \"24\u00d711=\",\n    \"14\u00d780=\",\n    \"34\u00d727=\",\n    \"44\u00d791=\",\n    \"47\u00d726=\",\n    \"82\u00d788=\",\n    \"96\u00d738=\",\n    \"36\u00d754=\",\n    \"35\u00d723=\",\n    \"16\u00d759=\",\n    \"46\u00d797=\",\n    \"39\u00d717=\",\n    \"19\u00d782=\",\n    \"78\u00d781=\",\n    \"76\u00d712=\",\n    \"31\u00d775=\",\n    \"63\u00d753=\",\n    \"46\u00d753=\",\n    \"99\u00d786=\",\n    \"58\u00d764=\",\n    \"97\u00d798=\",\n    \"43\u00d772=\",\n    \"100\u00d751=\",\n    \"43\u00d762=\",\n    \"16\u00d729=\",\n    \"96\u00d719=\",\n    \"84\u00d727=\",\n    \"72\u00d717=\"\n)\n\n$rowCount = $t.Rows.Count\n$colCount = $t.Columns.Count\n\nif (($rowCount * $colCount) -ne $oldVals.Count) {\n    throw \"Table has $($rowCount * $colCount) cells, expected $($oldVals.Count)\"\n}\n\n$i = 0\nfor ($r = 1; $r -le $rowCount; $r++) {\n    for ($c = 1; $c -le $colCount; $c++) {\n        $cell = $t.Cell($r, $c)\n        $cellRange = $cell.Range\n        # Trim the trailing end-of-cell marker(s) so Text reflects only\n        # the visible content for comparison.\n        $actual = $cellRange.Text.TrimEnd([char]13, [char]7)\n        $expectedOld = $oldVals[$i]\n        $newVal = $newVals[$i]\n\n        if ($actual -ne $expectedOld) {\n            throw \"Cell $i (row $r, col $c): expected '$expectedOld' but found '$actual'\"\n        }\n\n        $cellRange.Text = $newVal\n        $i++\n    }\n}\n"}
